$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- New row 7 data (fill before formulas so they evaluate correctly) ---
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 7
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

# --- Column C: 100-$D$11*Ax+$D$11 for rows 1-7 (each row gets its own formula text so it is not shared) ---
$ws.Range("C1").Formula = "=100-`$D`$11*A1+`$D`$11"
$ws.Range("C2").Formula = "=100-`$D`$11*A2+`$D`$11"
$ws.Range("C3").Formula = "=100-`$D`$11*A3+`$D`$11"
$ws.Range("C4").Formula = "=100-`$D`$11*A4+`$D`$11"
$ws.Range("C5").Formula = "=100-`$D`$11*A5+`$D`$11"
$ws.Range("C6").Formula = "=100-`$D`$11*A6+`$D`$11"
$ws.Range("C7").Formula = "=100-`$D`$11*A7+`$D`$11"

# --- Column F and I: shared formulas across rows 1-7 ---
$ws.Range("F1:F7").Formula = "=E1/D1"
$ws.Range("I1:I7").Formula = "=H1/G1"

# --- Column J ---
$ws.Range("J1").Formula = "=(C1*E1/D1)+(C1*H1/G1)"
$ws.Range("J2").Formula = "=`$J`$1*C2/100"
$ws.Range("J3").Formula = "=`$J`$1*C3/100"
$ws.Range("J4").Formula = "=`$J`$1*C4/100"
$ws.Range("J5").Formula = "=(`$J`$1*C5/100)"
$ws.Range("J6").Formula = "=`$J`$1*C6/100"
$ws.Range("J7").Formula = "=((`$J`$1*C7/100)*H7)"

# --- Column K: clear old text values in K1, then new formulas/values K2-K7 ---
$ws.Range("K1").ClearContents()
$ws.Range("K2").Formula = "=ROUND(J2*100/J1,0)"
$ws.Range("K3").Formula = "=ROUND(J3*100/J2,0)"
$ws.Range("K4").Formula = "=ROUND(J4*100/J3,0)"
$ws.Range("K5").Formula = "=ROUND(J5*100/J4,0)"
$ws.Range("K6").Formula = "=ROUND(J6*100/J5,0)"
$ws.Range("K7").Formula = "=ROUND(J7*100/J6,0)"

# --- Column L: new value ---
$ws.Range("L5").Value = 3

# --- Row 10 (new) ---
$ws.Range("K10").Value = 100

# --- Row 11 additions ---
$ws.Range("K11").Value = 66
$ws.Range("M11").Formula = "=K11*100/K10"

# --- Styling: add yellow highlight fill+font style for C3 and C5 ---
$ws.Range("C3").Interior.Color = 65535
$ws.Range("C3").Font.Color = 255
$ws.Range("C5").Interior.Color = 65535
$ws.Range("C5").Font.Color = 255

# --- Selection change ---
$ws.Range("L2").Select()

$wb.Save()
